# Updated cryptos list on Sat May 11 05:45:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the cell to be written as text even when the string looks like
    # a pure number (e.g. "584.96"), so Excel doesn't silently convert it
    # to a floating point number and lose the original formatting/precision.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.867.87"
$ws.Range("E2").Value = "  -3.02%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.917.88"
$ws.Range("E3").Value = "  -3.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "584.96"
$ws.Range("E5").Value = "  -1.53%  "

# Row 6 - Solana
Set-TextValue "D6" "145.27"
$ws.Range("E6").Value = "  -4.91%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.46%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.917.05"
$ws.Range("E9").Value = "  -3.60%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +4.30%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.144"
$ws.Range("E11").Value = "  -4.43%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -3.99%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.43%  "

# Row 14 - Avalanche
Set-TextValue "D14" "33.63"
$ws.Range("E14").Value = "  -5.50%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.10%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.401.40"
$ws.Range("E16").Value = "  -3.78%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "60.845.18"
$ws.Range("E17").Value = "  -3.10%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  -4.42%  "

# Row 19 - WrappedEther
Set-TextValue "D19" "2.918.80"
$ws.Range("E19").Value = "  -3.69%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "430.69"
$ws.Range("E20").Value = "  -4.85%  "

# Row 21 - Chainlink
Set-TextValue "D21" "13.61"
$ws.Range("E21").Value = "  -4.57%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.682"
$ws.Range("E22").Value = "  -2.19%  "

# Row 23 - Uniswap
Set-TextValue "D23" "7.13"
$ws.Range("E23").Value = "  -4.77%  "

# Row 24 - Litecoin
Set-TextValue "D24" "80.35"
$ws.Range("E24").Value = "  -3.31%  "

# Row 25 - RenderToken
Set-TextValue "D25" "10.84"
$ws.Range("E25").Value = "  -2.14%  "

# Row 26 - Fetch.AI
Set-TextValue "D26" "2.21"
$ws.Range("E26").Value = "  -3.47%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "11.93"
$ws.Range("E27").Value = "  -3.07%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.02%  "

# Row 29 - FirstDigitalUSD
$ws.Range("E29").Value = "  -0.02%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  -2.65%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -2.91%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "2.15"
$ws.Range("E32").Value = "  -3.56%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "26.55"
$ws.Range("E33").Value = "  -3.59%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -2.78%  "

# Row 35 - PEPE
Set-TextValue "D35" "0.0₃0870"
$ws.Range("E35").Value = "  +0.69%  "

# Row 36 - Mantle
$ws.Range("E36").Value = "  -2.58%  "

# Row 37 - Filecoin
Set-TextValue "D37" "5.65"
$ws.Range("E37").Value = "  -4.56%  "

# Row 38 and 39 swap coin identity (dogwifhat <-> Kaspa) with new values
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.129"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D39" "3.03"
$ws.Range("E39").Value = "  -4.32%  "

# Row 40 - OKB
Set-TextValue "D40" "49.69"
$ws.Range("E40").Value = "  -1.40%  "

# Row 41 - Stacks
Set-TextValue "D41" "2.00"
$ws.Range("E41").Value = "  -4.56%  "

# Row 42 - Cosmos
Set-TextValue "D42" "8.65"
$ws.Range("E42").Value = "  -4.95%  "

# Row 43 - TheGraph
Set-TextValue "D43" "0.296"
$ws.Range("E43").Value = "  -1.62%  "

# Row 44 - Arweave
Set-TextValue "D44" "40.72"
$ws.Range("E44").Value = "  -3.15%  "

# Row 45 - Bittensor
Set-TextValue "D45" "377.76"
$ws.Range("E45").Value = "  -4.26%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -2.58%  "

# Row 47 - Maker
Set-TextValue "D47" "2.678.86"
$ws.Range("E47").Value = "  -1.71%  "

# Row 48 - Monero
Set-TextValue "D48" "132.72"
$ws.Range("E48").Value = "  +0.75%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  -0.07%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "24.81"
$ws.Range("E50").Value = "  +1.68%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -1.76%  "
